$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.285.82"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "3.478.89"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.33%  "
$ws.Range("E7").Value = "  +4.09%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.91%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").Value = "4.036.09"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "3.487.21"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "65.473.66"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.988"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "86.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "609.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  +9.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0782"
$ws.Range("E39").Value = "  -6.77%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.364.11"
$ws.Range("E40").Value = "  +9.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.379"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  -8.38%  "
$ws.Range("E44").Value = "  -5.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0412"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "

Write-Host "Applied cryptos update."
